$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "(Version: 1.0.1)"
$ws.Range("A14").Value = "(Last tested with: ReportServer 4.0.0-6053) "
$ws.Range("A13").Select()
